# -----------------------------------------------------------------------
# Coinranking symbol-list refresh (GitHub Actions style update).
# Source data (columns D/E/G, plus the two swapped rows in B/C) are all
# stored as text in the workbook, so every write below targets cells we
# first force to Text (NumberFormat "@") -- this mirrors the upstream
# generator, which always emits these columns as inline strings, and
# keeps Excel from reinterpreting "301.38" / "-3.25%" / "17" as numbers.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng0 = $ws.Range("D2:D16")
$rng1 = $ws.Range("D19:D27")
$rng2 = $ws.Range("D39:D46")
$rng3 = $ws.Range("E2:E27")
$rng4 = $ws.Range("E39:E51")
$rng5 = $ws.Range("G2:G51")
$textTarget = $excel.Union($rng0, $rng1, $rng2, $rng3, $rng4, $rng5)
foreach ($area in $textTarget.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "301.38"
$ws.Range("E2").Value = "-3.25%"
$ws.Range("G2").Value = "17"
$ws.Range("D3").Value = "35.31"
$ws.Range("E3").Value = "-0.16%"
$ws.Range("G3").Value = "17"
$ws.Range("D4").Value = "5.060"
$ws.Range("E4").Value = "-0.47%"
$ws.Range("G4").Value = "17"
$ws.Range("D5").Value = "0.07964"
$ws.Range("E5").Value = "-2.37%"
$ws.Range("G5").Value = "17"
$ws.Range("D6").Value = "1.898"
$ws.Range("E6").Value = "-9.08%"
$ws.Range("G6").Value = "17"
$ws.Range("D7").Value = "7.779"
$ws.Range("E7").Value = "-2.19%"
$ws.Range("G7").Value = "17"
$ws.Range("D8").Value = "4.044"
$ws.Range("E8").Value = "-2.18%"
$ws.Range("G8").Value = "17"
$ws.Range("D9").Value = "0.9278"
$ws.Range("E9").Value = "-0.16%"
$ws.Range("G9").Value = "17"
$ws.Range("D10").Value = "0.1369"
$ws.Range("E10").Value = "32.11%"
$ws.Range("G10").Value = "17"
$ws.Range("D11").Value = "0.1903"
$ws.Range("E11").Value = "-0.29%"
$ws.Range("G11").Value = "17"
$ws.Range("D12").Value = "0.09067"
$ws.Range("E12").Value = "-1.26%"
$ws.Range("G12").Value = "17"
$ws.Range("D13").Value = "0.03429"
$ws.Range("E13").Value = "-6.15%"
$ws.Range("G13").Value = "17"
$ws.Range("D14").Value = "0.09915"
$ws.Range("E14").Value = "0.19%"
$ws.Range("G14").Value = "17"
$ws.Range("D15").Value = "0.001387"
$ws.Range("E15").Value = "-3.26%"
$ws.Range("G15").Value = "17"
$ws.Range("D16").Value = "0.005919"
$ws.Range("E16").Value = "1.76%"
$ws.Range("G16").Value = "17"
$ws.Range("E17").Value = "1.74%"
$ws.Range("G17").Value = "17"
$ws.Range("E18").Value = "-1.24%"
$ws.Range("G18").Value = "17"
$ws.Range("D19").Value = "0.3406"
$ws.Range("E19").Value = "-0.14%"
$ws.Range("G19").Value = "17"
$ws.Range("D20").Value = "0.1294"
$ws.Range("E20").Value = "-0.52%"
$ws.Range("G20").Value = "17"
$ws.Range("D21").Value = "5.075"
$ws.Range("E21").Value = "-0.83%"
$ws.Range("G21").Value = "17"
$ws.Range("D22").Value = "0.2397"
$ws.Range("E22").Value = "8.31%"
$ws.Range("G22").Value = "17"
$ws.Range("D23").Value = "0.04503"
$ws.Range("E23").Value = "-0.97%"
$ws.Range("G23").Value = "17"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "-1.19%"
$ws.Range("G24").Value = "17"
$ws.Range("D25").Value = "0.004771"
$ws.Range("E25").Value = "-0.02%"
$ws.Range("G25").Value = "17"
$ws.Range("D26").Value = "0.0001230"
$ws.Range("E26").Value = "-1.91%"
$ws.Range("G26").Value = "17"
$ws.Range("D27").Value = "0.0003000"
$ws.Range("E27").Value = "-32.76%"
$ws.Range("G27").Value = "17"
$ws.Range("G28").Value = "17"
$ws.Range("G29").Value = "17"
$ws.Range("G30").Value = "17"
$ws.Range("G31").Value = "17"
$ws.Range("G32").Value = "17"
$ws.Range("G33").Value = "17"
$ws.Range("G34").Value = "17"
$ws.Range("G35").Value = "17"
$ws.Range("G36").Value = "17"
$ws.Range("G37").Value = "17"
$ws.Range("G38").Value = "17"
$ws.Range("D39").Value = "0.01885"
$ws.Range("E39").Value = "-3.51%"
$ws.Range("G39").Value = "17"
$ws.Range("D40").Value = "0.04764"
$ws.Range("E40").Value = "-2.79%"
$ws.Range("G40").Value = "17"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.007282"
$ws.Range("E41").Value = "-3.37%"
$ws.Range("G41").Value = "17"
$ws.Range("B42").Value = "Dexo"
$ws.Range("C42").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").Value = "0.009616"
$ws.Range("E42").Value = "18.55%"
$ws.Range("G42").Value = "17"
$ws.Range("D43").Value = "0.1324"
$ws.Range("E43").Value = "-4.05%"
$ws.Range("G43").Value = "17"
$ws.Range("D44").Value = "0.002110"
$ws.Range("E44").Value = "-2.76%"
$ws.Range("G44").Value = "17"
$ws.Range("D45").Value = "0.01096"
$ws.Range("E45").Value = "-4.93%"
$ws.Range("G45").Value = "17"
$ws.Range("D46").Value = "0.00006229"
$ws.Range("E46").Value = "-5.70%"
$ws.Range("G46").Value = "17"
$ws.Range("E47").Value = "-0.32%"
$ws.Range("G47").Value = "17"
$ws.Range("E48").Value = "-65.08%"
$ws.Range("G48").Value = "17"
$ws.Range("E49").Value = "10.29%"
$ws.Range("G49").Value = "17"
$ws.Range("E50").Value = "-0.32%"
$ws.Range("G50").Value = "17"
$ws.Range("E51").Value = "-0.32%"
$ws.Range("G51").Value = "17"
